$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 89: remove the Reference (column N) entry, which pointed to the
#     same post id now duplicated in the new row. Fully clear cell (value + style)
#     so no empty <c> element remains in the OOXML. ---
$ws.Cells.Item(89, 14).Clear()

# --- Row 90: new entry for "DENSE FOG: NYE" post (31 Dec 2020) ---
$ws.Cells.Item(90, 1).Value = 44196
$ws.Cells.Item(90, 2).Value = 0.43055555555555558
$ws.Cells.Item(90, 3).Value = "Friends"
$ws.Cells.Item(90, 4).Value = "DENSE FOG: NYE ☁️ 🎆"
$ws.Cells.Item(90, 4).WrapText = $true
$ws.Cells.Item(90, 5).Value = "10107934723078259"
$ws.Cells.Item(90, 6).Value = 6
$ws.Cells.Item(90, 7).Value = 2
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 9).Value = 1
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 11).Value = 2
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 13).Value = 0

# --- Row 91: new entry for "DENSE FOG" post (2 Jan 2021) ---
$ws.Cells.Item(91, 1).Value = 44198
$ws.Cells.Item(91, 2).Value = 0.8208333333333333
$ws.Cells.Item(91, 3).Value = "Friends of Friends"
$ws.Cells.Item(91, 5).Value = "10107941020143889"
$ws.Cells.Item(91, 4).Value = "🌫 DENSE FOG 🌫"
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 1
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 1
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 13).Value = 1
$ws.Cells.Item(91, 15).Value = "Jeremy Ashton"

# --- Match final cursor/selection position left behind in the workbook ---
$ws.Range("I94").Select() | Out-Null
